# Clase 04/09/24 - ajustes finales
# - corrige tres valores de la tabla "Imagen con valores reemplazados en mediana"
#   (se repetian 232 por error y debian ser 226, igual que sus vecinos)
# - aplica formato numerico de 2 decimales a H35 (celda junto a la tabla de
#   dispersiones, dejada en blanco)
# - deja el cursor/seleccion en J8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = 226
$ws.Range("F19").Value = 226
$ws.Range("C21").Value = 226

$ws.Range("H35").NumberFormat = "0.00"

$ws.Range("J8").Select() | Out-Null
